$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[ Galit Zeluf%Andersson%null%1,                 Maria%Reinius%null%1,                 Lars E%Eriksson%null%1,                 Veronica%Svedhem%null%1,                 Farhad Mazi%Esfahani%null%1,                 Keshab%Deuba%null%1,                 Deepa%Rao%null%1,                 Goodluck Willey%Lyatuu%null%1,                 Danielle%Giovenco%null%1,                 Anna Mia%Ekström%null%1,                Galit Zeluf%Andersson%null%1,                Maria%Reinius%null%1,                Lars E%Eriksson%null%1,                Veronica%Svedhem%null%1,                Farhad Mazi%Esfahani%null%1,                Keshab%Deuba%null%1,                Deepa%Rao%null%1,                Goodluck Willey%Lyatuu%null%1,                Danielle%Giovenco%null%1,                Anna Mia%Ekström%null%1]"
$ws.Range("E3").Value = "[Roland Brian%Büchter%NULL%1,                Melanie%Messer%NULL%1]"
$ws.Range("E4").Value = "[Sarah%Clement%xref no email%1,          Francesca%Lassman%xref no email%1,          Elizabeth%Barley%xref no email%1,          Sara%Evans-Lacko%xref no email%1,          Paul%Williams%xref no email%1,          Sosei%Yamaguchi%xref no email%1,          Mike%Slade%xref no email%1,          Nicolas%R\u00fcsch%xref no email%1,          Graham%Thornicroft%xref no email%1]"
$ws.Range("E5").Value = "[Garumma Tolu%Feyissa%NULL%1,                Craig%Lockwood%NULL%2,                Craig%Lockwood%NULL%0,                Mirkuzie%Woldie%NULL%1,                Zachary%Munn%NULL%1,                Jerome A.%Singh%NULL%2,                Jerome A.%Singh%NULL%0,                NULL%NULL%NULL%0,                NULL%NULL%NULL%0]"
$ws.Range("E6").Value = "[Sabine E.%Hanisch%hanisch.sabine@siemens.com%1,                Conal D.%Twomey%c.twomey@soton.ac.uk%1,                Andrew C. H.%Szeto%aszeto@ucalgary.ca%1,                Ulrich W.%Birner%ulrich.birner@siemens.com%1,                Dennis%Nowak%dennis.nowak@med.uni-muenchen.de%1,                Carla%Sabariego%Carla.sabariego@med.uni-muenchen.de%1]"
$ws.Range("E7").Value = "[E.%Heim%NULL%1,                B. A.%Kohrt%NULL%2,                B. A.%Kohrt%NULL%0,                M.%Koschorke%NULL%1,                M.%Milenova%NULL%1,                G.%Thornicroft%NULL%1]"
$ws.Range("E8").Value = "[ Polly H. X.%Ma%null%1,        Zenobia C. Y.%Chan%null%1,        Alice Yuen%Loke%null%1]"
$ws.Range("E9").Value = "[ Harriet%Mills%null%1,                 Nadine%Mulfinger%null%1,                 Sophie%Raeder%null%1,                 Nicolas%Rüsch%null%1,                 Henry%Clements%null%1,                 Katrina%Scior%null%1,                Harriet%Mills%null%1,                Nadine%Mulfinger%null%1,                Sophie%Raeder%null%1,                Nicolas%Rüsch%null%1,                Henry%Clements%null%1,                Katrina%Scior%null%1]"
$ws.Range("E10").Value = "[Marija%Pantelic%NULL%1,                Janina I%Steinert%NULL%2,                Janina I%Steinert%NULL%0,                Jay%Park%NULL%1,                Shaun%Mellors%NULL%1,                Fungai%Murau%NULL%1]"
$ws.Range("E11").Value = "[Silatham%Sermrittirong%xref no email%1,          Wim H.%Van Brakel%xref no email%1,          J.F.G.%Bunbers-aelen%xref no email%1]"
